{"js": "// Remove the stray \"Prueba ariel\" run that was left over from a merge,\n// while keeping the (now empty) paragraph and its _GoBack bookmark intact.\n\nconst body = context.document.body;\n\n// Locate the run's text precisely so only that text is removed (the\n// paragraph mark, the bookmark, and the paragraph's own formatting stay).\nconst results = body.search(\"Prueba ariel\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the stray \"Prueba ariel\" run that was left over from a merge,\n# while keeping the (now empty) paragraph and its _GoBack bookmark intact.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Prueba ariel\"\n$find.Replacement.Text = \"\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceOne = 1 -> replace just the single match (the text of the run),\n# leaving the paragraph mark and bookmark untouched.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1) | Out-Null\n"}
